$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 45382
$ws.Range("C1").Value = 45291
$ws.Range("D1").Value = 45199
$ws.Range("E1").Value = 45107
$ws.Range("F1").Value = 45016

$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("F2").Value = 0

$ws.Range("B3").Value = 0.158
$ws.Range("C3").Value = 0.159
$ws.Range("D3").Value = 0.149715
$ws.Range("E3").Value = 0.125
$ws.Range("F3").Value = 0.149

$ws.Range("B4").Value = 30736000000
$ws.Range("C4").Value = 43221000000
$ws.Range("D4").Value = 30653000000
$ws.Range("E4").Value = 26783000000
$ws.Range("F4").Value = 31216000000

$ws.Range("B5").Value = 23636000000
$ws.Range("C5").Value = 33916000000
$ws.Range("D5").Value = 22956000000
$ws.Range("E5").Value = 19881000000
$ws.Range("F5").Value = 24160000000

$ws.Range("B6").Value = 2836000000
$ws.Range("C6").Value = 2848000000
$ws.Range("D6").Value = 2653000000
$ws.Range("E6").Value = 3052000000
$ws.Range("F6").Value = 2898000000

$ws.Range("B7").Value = 48482000000
$ws.Range("C7").Value = 64720000000
$ws.Range("D7").Value = 49071000000
$ws.Range("E7").Value = 45384000000
$ws.Range("F7").Value = 52860000000

$ws.Range("B8").Value = 30736000000
$ws.Range("C8").Value = 43221000000
$ws.Range("D8").Value = 30653000000
$ws.Range("E8").Value = 26783000000
$ws.Range("F8").Value = 31216000000

$ws.Range("B9").Value = 27900000000
$ws.Range("C9").Value = 40373000000
$ws.Range("D9").Value = 28000000000
$ws.Range("E9").Value = 23731000000
$ws.Range("F9").Value = 28318000000

$ws.Range("B10").Value = ""
$ws.Range("C10").Value = ""
$ws.Range("D10").Value = -18000000
$ws.Range("E10").Value = -18000000
$ws.Range("F10").Value = -12000000

$ws.Range("B11").Value = ""
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = 1002000000
$ws.Range("E11").Value = 998000000
$ws.Range("F11").Value = 930000000

$ws.Range("B12").Value = ""
$ws.Range("C12").Value = ""
$ws.Range("D12").Value = 984000000
$ws.Range("E12").Value = 980000000
$ws.Range("F12").Value = 918000000

$ws.Range("B13").Value = 23636000000
$ws.Range("C13").Value = 33916000000
$ws.Range("D13").Value = 22956000000
$ws.Range("E13").Value = 19881000000
$ws.Range("F13").Value = 24160000000

$ws.Range("B14").Value = 23636000000
$ws.Range("C14").Value = 33916000000
$ws.Range("D14").Value = 22956000000
$ws.Range("E14").Value = 19881000000
$ws.Range("F14").Value = 24160000000

$ws.Range("B15").Value = 62853000000
$ws.Range("C15").Value = 79202000000
$ws.Range("D15").Value = 62529000000
$ws.Range("E15").Value = 58799000000
$ws.Range("F15").Value = 66518000000

$ws.Range("B16").Value = 27900000000
$ws.Range("C16").Value = 40373000000
$ws.Range("D16").Value = 26969000000
$ws.Range("E16").Value = 22998000000
$ws.Range("F16").Value = 28318000000

$ws.Range("B17").Value = 15464709000
$ws.Range("C17").Value = 15576641000
$ws.Range("D17").Value = 15672400000
$ws.Range("E17").Value = 15775021000
$ws.Range("F17").Value = 15847050000

$ws.Range("B18").Value = 15405856000
$ws.Range("C18").Value = 15509763000
$ws.Range("D18").Value = 15599434000
$ws.Range("E18").Value = 15697614000
$ws.Range("F18").Value = 15787154000

$ws.Range("B19").Value = 1.53
$ws.Range("C19").Value = 2.18
$ws.Range("D19").Value = 1.46
$ws.Range("E19").Value = 1.26
$ws.Range("F19").Value = 1.52

$ws.Range("B20").Value = 1.53
$ws.Range("C20").Value = 2.19
$ws.Range("D20").Value = 1.47
$ws.Range("E20").Value = 1.27
$ws.Range("F20").Value = 1.53

$ws.Range("B21").Value = 23636000000
$ws.Range("C21").Value = 33916000000
$ws.Range("D21").Value = 22956000000
$ws.Range("E21").Value = 19881000000
$ws.Range("F21").Value = 24160000000

$ws.Range("B22").Value = 23636000000
$ws.Range("C22").Value = 33916000000
$ws.Range("D22").Value = 22956000000
$ws.Range("E22").Value = 19881000000
$ws.Range("F22").Value = 24160000000

$ws.Range("B23").Value = 23636000000
$ws.Range("C23").Value = 33916000000
$ws.Range("D23").Value = 22956000000
$ws.Range("E23").Value = 19881000000
$ws.Range("F23").Value = 24160000000

$ws.Range("B24").Value = 23636000000
$ws.Range("C24").Value = 33916000000
$ws.Range("D24").Value = 22956000000
$ws.Range("E24").Value = 19881000000
$ws.Range("F24").Value = 24160000000

$ws.Range("B25").Value = 23636000000
$ws.Range("C25").Value = 33916000000
$ws.Range("D25").Value = 22956000000
$ws.Range("E25").Value = 19881000000
$ws.Range("F25").Value = 24160000000

$ws.Range("B26").Value = 4422000000
$ws.Range("C26").Value = 6407000000
$ws.Range("D26").Value = 4042000000
$ws.Range("E26").Value = 2852000000
$ws.Range("F26").Value = 4222000000

$ws.Range("B27").Value = 28058000000
$ws.Range("C27").Value = 40323000000
$ws.Range("D27").Value = 26998000000
$ws.Range("E27").Value = 22733000000
$ws.Range("F27").Value = 28382000000

$ws.Range("B28").Value = 158000000
$ws.Range("C28").Value = -50000000
$ws.Range("D28").Value = 47000000
$ws.Range("E28").Value = -247000000
$ws.Range("F28").Value = 64000000

$ws.Range("B29").Value = 158000000
$ws.Range("C29").Value = -50000000
$ws.Range("D29").Value = 47000000
$ws.Range("E29").Value = -247000000
$ws.Range("F29").Value = 64000000

$ws.Range("B30").Value = ""
$ws.Range("C30").Value = ""
$ws.Range("D30").Value = -18000000
$ws.Range("E30").Value = -18000000
$ws.Range("F30").Value = -12000000

$ws.Range("B31").Value = ""
$ws.Range("C31").Value = ""
$ws.Range("D31").Value = 1002000000
$ws.Range("E31").Value = 998000000
$ws.Range("F31").Value = 930000000

$ws.Range("B32").Value = ""
$ws.Range("C32").Value = ""
$ws.Range("D32").Value = 984000000
$ws.Range("E32").Value = 980000000
$ws.Range("F32").Value = 918000000

$ws.Range("B33").Value = 27900000000
$ws.Range("C33").Value = 40373000000
$ws.Range("D33").Value = 26969000000
$ws.Range("E33").Value = 22998000000
$ws.Range("F33").Value = 28318000000

$ws.Range("B34").Value = 14371000000
$ws.Range("C34").Value = 14482000000
$ws.Range("D34").Value = 13458000000
$ws.Range("E34").Value = 13415000000
$ws.Range("F34").Value = 13658000000

$ws.Range("B35").Value = 7903000000
$ws.Range("C35").Value = 7696000000
$ws.Range("D35").Value = 7307000000
$ws.Range("E35").Value = 7442000000
$ws.Range("F35").Value = 7457000000

$ws.Range("B36").Value = 6468000000
$ws.Range("C36").Value = 6786000000
$ws.Range("D36").Value = 6151000000
$ws.Range("E36").Value = 5973000000
$ws.Range("F36").Value = 6201000000

$ws.Range("B37").Value = 42271000000
$ws.Range("C37").Value = 54855000000
$ws.Range("D37").Value = 40427000000
$ws.Range("E37").Value = 36413000000
$ws.Range("F37").Value = 41976000000

$ws.Range("B38").Value = 48482000000
$ws.Range("C38").Value = 64720000000
$ws.Range("D38").Value = 49071000000
$ws.Range("E38").Value = 45384000000
$ws.Range("F38").Value = 52860000000

$ws.Range("B39").Value = 90753000000
$ws.Range("C39").Value = 119575000000
$ws.Range("D39").Value = 89498000000
$ws.Range("E39").Value = 81797000000
$ws.Range("F39").Value = 94836000000

$ws.Range("B40").Value = 90753000000
$ws.Range("C40").Value = 119575000000
$ws.Range("D40").Value = 89498000000
$ws.Range("E40").Value = 81797000000
$ws.Range("F40").Value = 94836000000
